$wb = $excel.ActiveWorkbook

# --- Rename sheets (add suffix letters) ---
$wb.Worksheets.Item("AcessarUmProdutoPelaHome_P").Name = "AcessarUmProdutoPelaHome_Po"
$wb.Worksheets.Item("AcessarUmProdutoPelaHome_N").Name = "AcessarUmProdutoPelaHome_Ne"
$wb.Worksheets.Item("CadastrarNovoCliente_P").Name = "CadastrarNovoCliente_Po"
$wb.Worksheets.Item("CadastrarNovoCliente_N").Name = "CadastrarNovoCliente_Ne"
$wb.Worksheets.Item("BuscarUmProdutoPelaBusca_P").Name = "BuscarUmProdutoPelaBusca_Po"
$wb.Worksheets.Item("BuscarUmProdutoPelaBusca_N").Name = "BuscarUmProdutoPelaBusca_Ne"

# --- Sheet 1: AcessarUmProdutoPelaHome_Po ---
$ws1 = $wb.Worksheets.Item("AcessarUmProdutoPelaHome_Po")
$ws1.Range("A1").Value = "deveAbrirPaginaDeUmProdutoPelaCategoria"
$ws1.Range("B1").Value = "Categoria"
$ws1.Range("B2").Value = "HEADPHONES"
$ws1.Range("B3").Value = "LAPTOPS"
$ws1.Range("B4").Value = "SPEAKERS"
$ws1.Range("B5").Value = "TABLETS"
$ws1.Range("B6").Value = "MICE"

# --- Sheet 2: AcessarUmProdutoPelaHome_Ne ---
$ws2 = $wb.Worksheets.Item("AcessarUmProdutoPelaHome_Ne")
$ws2.Range("A1").Value = "deveTentarAbrirPaginaDeUmProdutoInesistentePelaCategoria"
$ws2.Range("B1").Value = "Categoria"
$ws2.Range("B2").Value = "HEADPHONES"
$ws2.Range("B3").Value = "LAPTOPS"
$ws2.Range("B4").Value = "SPEAKERS"
$ws2.Range("B5").Value = "TABLETS"
$ws2.Range("B6").Value = "MICE"

# --- Sheet 3: CadastrarNovoCliente_Po ---
$ws3 = $wb.Worksheets.Item("CadastrarNovoCliente_Po")
$ws3.Range("A1").Value = "DeveRealisarUmaBuscaComSucesso"

# --- Sheet 4: CadastrarNovoCliente_Ne ---
$ws4 = $wb.Worksheets.Item("CadastrarNovoCliente_Ne")
$ws4.Range("A1").Value = "DeveTentarCadasTrarUmNovoClienteComFalha"

# --- Sheet 6: BuscarUmProdutoPelaBusca_Ne ---
$ws6 = $wb.Worksheets.Item("BuscarUmProdutoPelaBusca_Ne")
$ws6.Range("A1").Value = "DeveRealisarUmaBuscaDeProdutoInvalido"

# --- Make sheet 4 (CadastrarNovoCliente_Ne) the active sheet/tab ---
$ws4.Activate()
$ws4.Range("C6:D6").Select()

$wb.Save()
